$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the date column as text so ambiguous dd/mm strings
# (e.g. "01/07/2025") are not auto-parsed into date serials.
$ws.Range("B2:B13").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value2 = 'Organisme : Direction regionale de l''agence nationale des eaux et forets de beni mellal-khenifra. Objet : La fourniture de la main d''oeuvre pour le gardiennage et la surveillance des forets relevant de la dranef de béni mellal -khénifra , réparti en quatre (04)lots.'
$ws.Cells.Item(2, 2).Value2 = '25/06/2025'
$ws.Cells.Item(2, 3).Value2 = 'N/A'

$ws.Cells.Item(3, 1).Value2 = 'Organisme : Centre hospitalier préfectoral de m’diq-fnideq. Objet : L’hygiène, et le nettoyage des locaux et des services du centre hospitalier préfectoral m’diq fnideq. (hopital mohammed vi m’diq et hopital hassan ii fnideq)'
$ws.Cells.Item(3, 2).Value2 = '25/06/2025'
$ws.Cells.Item(3, 3).Value2 = 'N/A'

$ws.Cells.Item(4, 1).Value2 = 'Organisme : Direction regionale de l’agriculture marrakech safi. Objet : Appel d''offres ouvert international à majoration, pour : gardiennage, surveillance et sécurité des locaux des directions provinciales de l’agriculture de marrakech, rhamna, chichaoua, essaouira, safi et la direction régionale de l’agriculture marrakech safi.'
$ws.Cells.Item(4, 2).Value2 = '30/06/2025'
$ws.Cells.Item(4, 3).Value2 = 'N/A'

$ws.Cells.Item(5, 1).Value2 = 'Organisme : Societe tanger med utilities. Objet : Gardiennage et surveillance des chantiers aux zones gérées par tmu tanger med utilities (tmu).'
$ws.Cells.Item(5, 2).Value2 = '01/07/2025'
$ws.Cells.Item(5, 3).Value2 = 'N/A'

$ws.Cells.Item(6, 1).Value2 = 'Organisme : Direction du port de casablanca et region. Objet : Prestation de gardiennage et surveillance du port de casablanca par les maitres-chiens'
$ws.Cells.Item(6, 2).Value2 = '01/07/2025'
$ws.Cells.Item(6, 3).Value2 = 'N/A'

$ws.Cells.Item(7, 1).Value2 = 'Organisme : Direction provinciale de tiznit. Objet : Appel d’offres ouvert international a majoration pour prestations de surveillance et de gardiennage des installations sportives a la province de tiznit en 4 lots:'
$ws.Cells.Item(7, 2).Value2 = '03/07/2025'
$ws.Cells.Item(7, 3).Value2 = 'N/A'

$ws.Cells.Item(8, 1).Value2 = 'Organisme : Centre hospitalier provincial de kenitra. Objet : Appel d''offres ouvert a majoration : activités d’accueil des malades du centre hospitalier provincial de kenitra.'
$ws.Cells.Item(8, 2).Value2 = '07/07/2025'
$ws.Cells.Item(8, 3).Value2 = 'N/A'

$ws.Cells.Item(9, 1).Value2 = 'Organisme : Centre hospitalier provincial de tetouan. Objet : Execution des prestations de gardiennage des batiments relevant du centre hospitalier provincial de tetouan (hopital civil - hopital ben karrich et hopital errazi de tetouan)'
$ws.Cells.Item(9, 2).Value2 = '08/07/2025'
$ws.Cells.Item(9, 3).Value2 = 'N/A'

$ws.Cells.Item(10, 1).Value2 = 'Organisme : Centre hospitalier regional d''agadir. Objet : Les presttions des activites d''accueil pour le centre hospitalier régional d''agadir.'
$ws.Cells.Item(10, 2).Value2 = '22/07/2025'
$ws.Cells.Item(10, 3).Value2 = 'N/A'

$ws.Cells.Item(11, 1).Value2 = 'Organisme : Direction des affaires administratives. Objet : Appel d’offres ouvert international, à majoration n° 16/2025, pour l’entretien et nettoyage des batiments et locaux administratifs du service central du ministère de l’intérieur et ses annexes à rabat et salé (lot unique).'
$ws.Cells.Item(11, 2).Value2 = '28/07/2025'
$ws.Cells.Item(11, 3).Value2 = 'N/A'

$ws.Cells.Item(12, 1).Value2 = 'Organisme : Dr de l''artisanat region rabat-zemour-zaier. Objet : Ao ouvert international a majoration : prestations de gardiennage et de surveillance des locaux de la direction régionale de l’artisanat et de l’économie sociale et solidaire de rabat et des etablissements de la formation professionnelles et des agences relevant du secrétariat d’etat chargé de l’artisanat et de l’economie sociale et solidaire'
$ws.Cells.Item(12, 2).Value2 = '28/07/2025'
$ws.Cells.Item(12, 3).Value2 = 'N/A'

$ws.Cells.Item(13, 1).Value2 = 'Organisme : Direction régionale de l’artisanat d’agadir. Objet : Ao ouvert international a majoration : gardiennage et surveillance des batiments administratifs relevant de la direction régionale de l’artisanat et de l’economie sociale souss massa et ses entités en lot unique.'
$ws.Cells.Item(13, 2).Value2 = '05/08/2025'
$ws.Cells.Item(13, 3).Value2 = 'N/A'

# Restore the default (unstyled) look for the date column.
$ws.Range("B2:B13").Style = "Normal"